$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Save" column header in H1 - clone the formatting of the neighboring
# header cell (G1: bold, centered, bordered) so the new header matches the
# existing header row style exactly.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Row 2 data value for the new Save column.
$ws.Range("H2").Value = 1
